$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 270
$ws.Range("C3").Value = 167794
$ws.Range("C4").Value = 158674
$ws.Range("C8").Value = 65.3
